# Adds the "total_release_deploy_raw" column (K) to both the Octubre and
# Noviembre sheets, and refreshes several existing raw metric values that
# moved as a result of recomputing the underlying report.

function Set-SheetData($ws, $kValues, $changes) {
    # New column header
    $ws.Range("K1").Value = "total_release_deploy_raw"

    # New column K values (rows 2..14)
    for ($i = 0; $i -lt $kValues.Length; $i++) {
        $ws.Cells.Item($i + 2, 11).Value = $kValues[$i]
    }

    # Refreshed values in existing columns
    foreach ($change in $changes) {
        $ws.Cells.Item($change[0], $change[1]).Value = $change[2]
    }
}

$wb = $excel.ActiveWorkbook
$wsOctubre = $wb.Worksheets.Item("Octubre")
$wsNoviembre = $wb.Worksheets.Item("Noviembre")

$octubreK = @(0, 7, 0, 6, 10, 5, 27, 5, 3, 0, 0, 0, 0)

$octubreChanges = @(
    ,@(7, 7, 23.75)
    ,@(7, 8, 2.72)
    ,@(8, 2, 17667)
    ,@(8, 6, 13112)
    ,@(9, 2, 12031)
    ,@(9, 6, 9465)
    ,@(11, 2, 10127)
    ,@(11, 6, 7262)
    ,@(11, 7, 12.11)
    ,@(11, 8, 8.2)
    ,@(13, 2, 6761)
    ,@(13, 6, 4732)
)

$noviembreK = @(0, 9, 7, 2, 11, 5, 18, 2, 3, 0, 0, 0, 0)

$noviembreChanges = @(
    ,@(2, 2, 2230)
    ,@(2, 3, 1184)
    ,@(2, 4, 145)
    ,@(2, 5, 901)
    ,@(2, 7, 3.51)
    ,@(2, 8, 15.47)
    ,@(3, 2, 9496)
    ,@(3, 3, 713)
    ,@(3, 4, 241)
    ,@(3, 5, 728)
    ,@(3, 7, 9.59)
    ,@(3, 8, 22.32)
    ,@(4, 2, 13292)
    ,@(4, 3, 359)
    ,@(4, 4, 231)
    ,@(4, 5, 158)
    ,@(4, 6, 12544)
    ,@(4, 7, 14.75)
    ,@(4, 8, 6.3)
    ,@(5, 2, 22623)
    ,@(5, 3, 1440)
    ,@(5, 4, 57)
    ,@(5, 5, 592)
    ,@(5, 6, 20534)
    ,@(5, 7, 15.63)
    ,@(5, 8, 5.03)
    ,@(6, 2, 6889)
    ,@(6, 3, 1401)
    ,@(6, 4, 187)
    ,@(6, 5, 644)
    ,@(6, 7, 7.41)
    ,@(6, 8, 6.45)
    ,@(7, 2, 16710)
    ,@(7, 3, 1561)
    ,@(7, 4, 442)
    ,@(7, 5, 2267)
    ,@(7, 7, 21.45)
    ,@(7, 8, 5.41)
    ,@(8, 2, 9668)
    ,@(8, 3, 619)
    ,@(8, 4, 360)
    ,@(8, 5, 1412)
    ,@(8, 6, 7277)
    ,@(8, 7, 11.11)
    ,@(8, 8, 4.32)
    ,@(9, 2, 7934)
    ,@(9, 3, 2615)
    ,@(9, 4, 314)
    ,@(9, 5, 867)
    ,@(9, 7, 9.06)
    ,@(9, 8, 6.66)
    ,@(10, 2, 16183)
    ,@(10, 3, 241)
    ,@(10, 4, 32)
    ,@(10, 5, 223)
    ,@(10, 6, 15687)
    ,@(10, 7, 4.22)
    ,@(10, 8, 18.65)
    ,@(11, 2, 13223)
    ,@(11, 3, 846)
    ,@(11, 4, 344)
    ,@(11, 5, 1279)
    ,@(11, 6, 10753)
    ,@(11, 7, 15.77)
    ,@(11, 8, 5.34)
    ,@(12, 2, 15231)
    ,@(12, 3, 1027)
    ,@(12, 4, 92)
    ,@(12, 5, 486)
    ,@(12, 6, 13626)
    ,@(12, 7, 9.09)
    ,@(12, 8, 10.04)
    ,@(13, 2, 5082)
    ,@(13, 3, 1899)
    ,@(13, 4, 229)
    ,@(13, 7, 6.28)
    ,@(13, 8, 11.06)
    ,@(14, 2, 9496)
    ,@(14, 3, 713)
    ,@(14, 4, 241)
    ,@(14, 5, 728)
    ,@(14, 7, 9.59)
    ,@(14, 8, 22.32)
)

Set-SheetData $wsOctubre $octubreK $octubreChanges
Set-SheetData $wsNoviembre $noviembreK $noviembreChanges
